# Daily attendance processing - 2025-11-05 04:49:17
# Normalizes the "Recorded By" column (G): the most recent recorder
# (previously appended last) is moved to the front of the comma-separated
# list so the latest recorder is listed first.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $val = $cell.Value2

    if ($val -ne $null -and $val -ne "") {
        $parts = $val -split ", "
        if ($parts.Count -gt 1) {
            $last = $parts[$parts.Count - 1]
            $rest = $parts[0..($parts.Count - 2)]
            $newVal = $last + ", " + ($rest -join ", ")
            $cell.Value2 = $newVal
        }
    }
}
